# Change table declarations from id='...' to class='...'
# e.g. !!ObjTables type='Data' id='Test'  ->  !!ObjTables type='Data' class='Test'

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string] -and $val -like "*!!ObjTables type='Data' id='*") {
            $cell.Value2 = $val -replace "!!ObjTables type='Data' id='", "!!ObjTables type='Data' class='"
        }
    }
}
